# Weekly update: insert a new price record for "Ají" at
# Terminal Hortofrutícola Agro Chillán.
#
# The new observation is inserted as row 35, pushing the existing rows
# 35-43 down to 36-44 (a classic "insert row" at position 35).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 35; everything below (old rows 35-43) shifts down
# to 36-44, and the sheet's used range grows to R44.
$ws.Rows(35).Insert()

# Populate the newly-inserted row 35 with this week's data.
$ws.Range("A35").Value = 7
$ws.Range("B35").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C35").Value = "Ñuble"
$ws.Range("D35").Value = 44543
$ws.Range("E35").Value = 16
$ws.Range("F35").Value = 100112021
$ws.Range("G35").Value = "Ají"
$ws.Range("H35").Value = "Americana (o)"
$ws.Range("I35").Value = "Primera"
$ws.Range("J35").Value = 100
$ws.Range("K35").Value = 15000
$ws.Range("L35").Value = 16000
$ws.Range("M35").Value = 15500
$ws.Range("N35").Value = "$/caja 15 kilos"
$ws.Range("O35").Value = "Región del Maule"
$ws.Range("P35").Value = 1033
$ws.Range("Q35").Value = 15
$ws.Range("R35").Value = "Hortaliza"
